$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.371340870857239
$ws.Range("B1").Value = 1.630016088485718
$ws.Range("C1").Value = 3.163380146026611
$ws.Range("D1").Value = 1.51623547077179
$ws.Range("E1").Value = 0.8312529325485229
